# Update the EPEX Spot price workbook with the newest daily data point.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new date column (14-sep) after the last one (CN).
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("CO1").Value = "14-sep"

# Copy the header formatting (bold font, borders, centered) from CN1 to CO1.
$wsSpot.Range("CN1").Copy()
$wsSpot.Range("CO1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$spotValues = @{
    2  = 45.73
    3  = 22.17
    4  = 17.96
    5  = 19.06
    6  = 14.42
    7  = 16.74
    8  = 19.38
    9  = 19.74
    10 = 18.9
    11 = 23.58
    12 = 15.43
    13 = 17.33
    14 = 6.5
    15 = 0
    16 = -0.01
    17 = 0
    18 = 5.59
    19 = 4.56
    20 = 15.65
    21 = 18.4
    22 = 16.79
    23 = 13.53
    24 = 17.36
    25 = 12.92
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item($row, 93).Value = $spotValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the newest daily gas price row.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force the date to be stored as plain text (matching the existing column),
# not auto-converted to a date serial number, then drop the temporary
# number-format override so the cell keeps the default (unstyled) look.
$wsGaz.Range("A90").NumberFormat = "@"
$wsGaz.Range("A90").Value = "2025-09-12"
$wsGaz.Range("A90").ClearFormats()
$wsGaz.Range("B90").Value = 32.2

# ---------------------------------------------------------------------------
# Sheet "CO2": append the newest daily CO2 price row.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A90").NumberFormat = "@"
$wsCo2.Range("A90").Value = "2025-09-12"
$wsCo2.Range("A90").ClearFormats()
$wsCo2.Range("B90").Value = 75.47
